$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.315.84'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.880.73'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.59'
$ws.Range('E5').Value = '  -2.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.681'
$ws.Range('E6').Value = '  -2.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.69'
$ws.Range('E8').Value = '  +4.89%  '
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '53.29'
$ws.Range('E10').Value = '  +1.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0742'
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.49'
$ws.Range('E13').Value = '  +3.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.152.98'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.766'
$ws.Range('E15').Value = '  +4.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.94'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.865.41'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.300.85'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.53'
$ws.Range('E19').Value = '  -0.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0824'
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.61'
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.04'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('E24').Value = '  +11.01%  '
$ws.Range('E25').Value = '  -0.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -5.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.90'
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.63'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.84'
$ws.Range('E35').Value = '  -11.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.43'
$ws.Range('E36').Value = '  -12.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.854'
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('E38').Value = '  -4.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0726'
$ws.Range('E39').Value = '  +9.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.39'
$ws.Range('E40').Value = '  -1.39%  '
$ws.Range('E41').Value = '  +1.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.83'
$ws.Range('E42').Value = '  -1.99%  '
$ws.Range('E43').Value = '  -3.10%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.41'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.306.73'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0798'
$ws.Range('E46').Value = '  +4.59%  '
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.73'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.84'
$ws.Range('E49').Value = '  -2.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.30'
$ws.Range('E50').Value = '  -4.76%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.054.93'
$ws.Range('E51').Value = '  -1.44%  '
